$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "97.318.23"
$ws.Range("E2").Value = "  +2.16%  "

# Row 3
$ws.Range("D3").Value = "3.583.39"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.15"
$ws.Range("E5").Value = "  +2.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "658.59"
$ws.Range("E6").Value = "  +1.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.70"
$ws.Range("E7").Value = "  +15.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.427"
$ws.Range("E8").Value = "  +6.82%  "

# Row 9
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.05"
$ws.Range("E10").Value = "  +4.30%  "

# Row 11
$ws.Range("D11").Value = "3.580.01"
$ws.Range("E11").Value = "  +0.05%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.58"
$ws.Range("E12").Value = "  +5.11%  "

# Row 13
$ws.Range("E13").Value = "  +0.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  -0.33%  "

# Row 15
$ws.Range("D15").Value = "4.250.71"
$ws.Range("E15").Value = "  +0.13%  "

# Row 16
$ws.Range("D16").Value = "97.075.56"
$ws.Range("E16").Value = "  +2.05%  "

# Row 17
$ws.Range("E17").Value = "  +3.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.66"
$ws.Range("E18").Value = "  +9.92%  "

# Row 19
$ws.Range("D19").Value = "3.584.36"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20
$ws.Range("E20").Value = "  +1.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.05"
$ws.Range("E21").Value = "  +1.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.528"
$ws.Range("E22").Value = "  +10.24%  "

# Row 23
$ws.Range("E23").Value = "  +1.31%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "516.37"
$ws.Range("E24").Value = "  +1.54%  "

# Row 25
$ws.Range("E25").Value = "  +4.85%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("E26").Value = "  +0.79%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.10"
$ws.Range("E27").Value = "  +6.05%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.09"
$ws.Range("E28").Value = "  +3.94%  "

# Row 29
$ws.Range("D29").Value = "3.775.26"
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.164"
$ws.Range("E30").Value = "  +14.76%  "

# Row 31
$ws.Range("E31").Value = "  -0.69%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.90"
$ws.Range("E32").Value = "  +3.74%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.996"
$ws.Range("E33").Value = "  -0.67%  "

# Row 34
$ws.Range("E34").Value = "  +3.84%  "

# Row 35
$ws.Range("E35").Value = "  +0.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.74"
$ws.Range("E36").Value = "  -0.27%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.90"
$ws.Range("E37").Value = "  +3.95%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "620.82"
$ws.Range("E38").Value = "  +6.40%  "

# Row 39
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.568"
$ws.Range("E39").Value = "  +1.35%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.66"
$ws.Range("E40").Value = "  -0.69%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.154"
$ws.Range("E41").Value = "  +2.57%  "

# Row 42
$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.95"
$ws.Range("E42").Value = "  +8.79%  "

# Row 43
$ws.Range("E43").Value = "  -0.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.923"
$ws.Range("E44").Value = "  +2.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.00"
$ws.Range("E45").Value = "  +4.41%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0440"
$ws.Range("E46").Value = "  +5.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("E48").Value = "  +0.94%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.408"
$ws.Range("E49").Value = "  +30.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.51"
$ws.Range("E50").Value = "  +4.56%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.31"
$ws.Range("E51").Value = "  +8.06%  "
